$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new text value to cell E6
$ws.Range("E6").Value = "Bla bla bla bla bla"

# Update the active selection to match the final state (E7)
$ws.Range("E7").Select()
